# edit.ps1 — applies the scheduled-runner profit recalculation to Shinryu_Profits.xlsx
# Updates currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N) on 29 leve rows
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets, per the refreshed market-board pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 223.875  # row 2: Mercury Rising / Quicksilver
$ws.Range("I2").Value = 223.875
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 223.875
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = -110.875

$ws.Range("H33").Value = 76.083336  # row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("I33").Value = 77.14286
$ws.Range("J33").Value = 74.59999999999999
$ws.Range("K33").Value = 77.14286
$ws.Range("L33").Value = 74.59999999999999
$ws.Range("M33").Value = 151.85714
$ws.Range("N33").Value = -532.6

$ws.Range("H40").Value = 1181.2222  # row 40: Stuck in the Moment / Horn Glue
$ws.Range("I40").Value = 1200.1666
$ws.Range("J40").Value = 1143.3334
$ws.Range("K40").Value = 1200.1666
$ws.Range("L40").Value = 1143.3334
$ws.Range("M40").Value = -1025.1666
$ws.Range("N40").Value = -1493.3334

$ws.Range("H69").Value = 24392628  # row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("I69").Value = 2350.4
$ws.Range("K69").Value = 7051.200000000001
$ws.Range("M69").Value = -6177.200000000001

$ws.Range("H72").Value = 24392628  # row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("I72").Value = 2350.4
$ws.Range("K72").Value = 21153.6
$ws.Range("M72").Value = -16785.6

$ws.Range("H86").Value = 13335333  # row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("I86").Value = 1833.3334
$ws.Range("J86").Value = 66669332
$ws.Range("K86").Value = 1833.3334
$ws.Range("L86").Value = 66669332
$ws.Range("M86").Value = -710.3334
$ws.Range("N86").Value = -66671578

$ws.Range("H89").Value = 13335333  # row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("I89").Value = 1833.3334
$ws.Range("J89").Value = 66669332
$ws.Range("K89").Value = 9166.666999999999
$ws.Range("L89").Value = 333346660
$ws.Range("M89").Value = -3550.666999999999
$ws.Range("N89").Value = -333357892

$ws.Range("H138").Value = 2896.4255  # row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("I138").Value = 1231.1578
$ws.Range("J138").Value = 4026.4285
$ws.Range("K138").Value = 3693.4734
$ws.Range("L138").Value = 12079.2855
$ws.Range("M138").Value = 1446.5266
$ws.Range("N138").Value = -22359.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2762.3333  # row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("I61").Value = 2641.1667
$ws.Range("J61").Value = 3004.6667
$ws.Range("K61").Value = 2641.1667
$ws.Range("L61").Value = 3004.6667
$ws.Range("M61").Value = -2429.1667
$ws.Range("N61").Value = -3428.6667

$ws.Range("H136").Value = 2762.3333  # row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("I136").Value = 2641.1667
$ws.Range("J136").Value = 3004.6667
$ws.Range("K136").Value = 7923.500100000001
$ws.Range("L136").Value = 9014.000100000001
$ws.Range("M136").Value = -5373.500100000001
$ws.Range("N136").Value = -14114.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 45000  # row 125: Archon of His Eye / High Durium Knives
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 776.6818  # row 68: Such a Butter Face / Fermented Butter
$ws.Range("I68").Value = 708.6389
$ws.Range("J68").Value = 858.3333
$ws.Range("K68").Value = 2125.9167
$ws.Range("L68").Value = 2574.9999
$ws.Range("M68").Value = -1314.9167
$ws.Range("N68").Value = -4196.9999

$ws.Range("H71").Value = 776.6818  # row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("I71").Value = 708.6389
$ws.Range("J71").Value = 858.3333
$ws.Range("K71").Value = 6377.7501
$ws.Range("L71").Value = 7724.9997
$ws.Range("M71").Value = -2321.7501
$ws.Range("N71").Value = -15836.9997

$ws.Range("H131").Value = 2468.7424  # row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("I131").Value = 769
$ws.Range("J131").Value = 2578.4033
$ws.Range("K131").Value = 2307
$ws.Range("L131").Value = 7735.2099
$ws.Range("M131").Value = 2733
$ws.Range("N131").Value = -17815.2099

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 5000  # row 19: Better Four Eyes than None / Brass Spectacles
$ws.Range("I19").Value = 5000
$ws.Range("K19").Value = 5000
$ws.Range("M19").Value = -4712

$ws.Range("H80").Value = 2974.3914  # row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("J80").Value = 3414
$ws.Range("L80").Value = 3414
$ws.Range("N80").Value = -5410

$ws.Range("H83").Value = 2974.3914  # row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("J83").Value = 3414
$ws.Range("L83").Value = 17070
$ws.Range("N83").Value = -27054

$ws.Range("H102").Value = 1629.75  # row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("I102").Value = 1560.3334
$ws.Range("J102").Value = 1838
$ws.Range("K102").Value = 1560.3334
$ws.Range("L102").Value = 1838
$ws.Range("M102").Value = 61.66660000000002
$ws.Range("N102").Value = -5082

$ws.Range("H132").Value = 3438.9678  # row 132: On Board for Lar / Lar Ingot
$ws.Range("I132").Value = 3322.5557
$ws.Range("J132").Value = 4224.75
$ws.Range("K132").Value = 9967.667099999999
$ws.Range("L132").Value = 12674.25
$ws.Range("M132").Value = -7437.667099999999
$ws.Range("N132").Value = -17734.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 66667644  # row 46: Supply Side Logic / Boar Leather
$ws.Range("I46").Value = 76924056
$ws.Range("K46").Value = 76924056
$ws.Range("M46").Value = -76923868

$ws.Range("H55").Value = 296.5909  # row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("I55").Value = 207.14285
$ws.Range("J55").Value = 338.33334
$ws.Range("K55").Value = 207.14285
$ws.Range("L55").Value = 338.33334
$ws.Range("M55").Value = -34.14285000000001
$ws.Range("N55").Value = -684.33334

$ws.Range("H68").Value = 418405.5  # row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("J68").Value = 1743.3334
$ws.Range("L68").Value = 1743.3334
$ws.Range("N68").Value = -3241.3334

$ws.Range("H71").Value = 418405.5  # row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("J71").Value = 1743.3334
$ws.Range("L71").Value = 8716.666999999999
$ws.Range("N71").Value = -16204.667

$ws.Range("H82").Value = 2185.7646  # row 82: Trainin' the Neck / Dragon Leather
$ws.Range("I82").Value = 2329.8
$ws.Range("J82").Value = 1980
$ws.Range("K82").Value = 2329.8
$ws.Range("L82").Value = 1980
$ws.Range("M82").Value = -1968.8
$ws.Range("N82").Value = -2702

$ws.Range("H85").Value = 2185.7646  # row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("I85").Value = 2329.8
$ws.Range("J85").Value = 1980
$ws.Range("K85").Value = 2329.8
$ws.Range("L85").Value = 1980
$ws.Range("M85").Value = -1081.8
$ws.Range("N85").Value = -4476

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 9999.5  # row 31: Whatchoo Talking About / Cotton Doublet Vest of Crafting
$ws.Range("J31").Value = 9999.5
$ws.Range("L31").Value = 9999.5
$ws.Range("N31").Value = -10695.5

$ws.Range("H81").Value = 1161.3334  # row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("I81").Value = 1180.2
$ws.Range("J81").Value = 1137.75
$ws.Range("K81").Value = 2360.4
$ws.Range("L81").Value = 2275.5
$ws.Range("M81").Value = -1299.4
$ws.Range("N81").Value = -4397.5

$ws.Range("H84").Value = 1161.3334  # row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("I84").Value = 1180.2
$ws.Range("J84").Value = 1137.75
$ws.Range("K84").Value = 11802
$ws.Range("L84").Value = 11377.5
$ws.Range("M84").Value = -6498
$ws.Range("N84").Value = -21985.5

$ws.Range("H126").Value = 1535.8182  # row 126: A Polished Purchase / Snow Linen
$ws.Range("I126").Value = 1169.1666
$ws.Range("J126").Value = 1975.8
$ws.Range("K126").Value = 3507.4998
$ws.Range("L126").Value = 5927.4
$ws.Range("M126").Value = -1037.4998
$ws.Range("N126").Value = -10867.4
